$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2195121951219512
$ws.Range("C2").Value = 0.5284552845528455
$ws.Range("J2").Value = 0.01626016260162602
$ws.Range("O2").Value = 0.004065040650406504
$ws.Range("P2").Value = 0.1382113821138211
$ws.Range("S2").Value = 0.09349593495934959
$ws.Range("C3").Value = 0.02272727272727273
$ws.Range("J3").Value = 0.02272727272727273
$ws.Range("P3").Value = 0.7803030303030303
$ws.Range("S3").Value = 0.1742424242424243
$ws.Range("J4").Value = 0.02857142857142857
$ws.Range("P4").Value = 0.6571428571428571
$ws.Range("S4").Value = 0.3142857142857143
$ws.Range("B6").Value = 0.03937007874015748
$ws.Range("D6").Value = 0.01181102362204724
$ws.Range("F6").Value = 0.08267716535433071
$ws.Range("J6").Value = 0.2283464566929134
$ws.Range("O6").Value = 0.01574803149606299
$ws.Range("Q6").Value = 0.1653543307086614
$ws.Range("R6").Value = 0.07874015748031496
$ws.Range("S6").Value = 0.3779527559055118
$ws.Range("B7").Value = 0.1105263157894737
$ws.Range("D7").Value = 0.02105263157894737
$ws.Range("F7").Value = 0.08421052631578947
$ws.Range("J7").Value = 0.08421052631578947
$ws.Range("O7").Value = 0.01578947368421053
$ws.Range("Q7").Value = 0.1368421052631579
$ws.Range("R7").Value = 0.1421052631578947
$ws.Range("S7").Value = 0.4052631578947368
$ws.Range("B8").Value = 0.08695652173913043
$ws.Range("D8").Value = 0.007561436672967864
$ws.Range("E8").Value = 0.005671077504725898
$ws.Range("F8").Value = 0.05671077504725898
$ws.Range("J8").Value = 0.1077504725897921
$ws.Range("O8").Value = 0.007561436672967864
$ws.Range("Q8").Value = 0.1776937618147448
$ws.Range("R8").Value = 0.1077504725897921
$ws.Range("S8").Value = 0.44234404536862
$ws.Range("B9").Value = 0.04736842105263158
$ws.Range("D9").Value = 0.03157894736842105
$ws.Range("F9").Value = 0.05789473684210526
$ws.Range("J9").Value = 0.1
$ws.Range("O9").Value = 0.02631578947368421
$ws.Range("Q9").Value = 0.1947368421052632
$ws.Range("R9").Value = 0.09473684210526316
$ws.Range("S9").Value = 0.4473684210526316
$ws.Range("B10").Value = 0.08847402597402597
$ws.Range("D10").Value = 0.01623376623376623
$ws.Range("E10").Value = 0.002435064935064935
$ws.Range("F10").Value = 0.08603896103896104
$ws.Range("J10").Value = 0.09577922077922078
$ws.Range("O10").Value = 0.0211038961038961
$ws.Range("Q10").Value = 0.213474025974026
$ws.Range("R10").Value = 0.09983766233766234
$ws.Range("S10").Value = 0.3766233766233766
$ws.Range("G11").Value = 0.1278688524590164
$ws.Range("J11").Value = 0.1311475409836066
$ws.Range("K11").Value = 0.2
$ws.Range("L11").Value = 0.5245901639344263
$ws.Range("S11").Value = 0.01639344262295082
$ws.Range("G12").Value = 0.725609756097561
$ws.Range("J12").Value = 0.2317073170731707
$ws.Range("K12").Value = 0.01219512195121951
$ws.Range("L12").Value = 0.02439024390243903
$ws.Range("S12").Value = 0.006097560975609756
$ws.Range("G13").Value = 0.7804878048780488
$ws.Range("J13").Value = 0.2195121951219512
$ws.Range("F15").Value = 0.02403846153846154
$ws.Range("H15").Value = 0.2115384615384615
$ws.Range("I15").Value = 0.07692307692307693
$ws.Range("J15").Value = 0.3413461538461539
$ws.Range("K15").Value = 0.0625
$ws.Range("M15").Value = 0.01442307692307692
$ws.Range("O15").Value = 0.02403846153846154
$ws.Range("S15").Value = 0.2451923076923077
$ws.Range("F16").Value = 0.02531645569620253
$ws.Range("H16").Value = 0.1518987341772152
$ws.Range("I16").Value = 0.08227848101265822
$ws.Range("J16").Value = 0.4050632911392405
$ws.Range("K16").Value = 0.120253164556962
$ws.Range("M16").Value = 0.0189873417721519
$ws.Range("N16").Value = 0.006329113924050633
$ws.Range("O16").Value = 0.05696202531645569
$ws.Range("S16").Value = 0.1329113924050633
$ws.Range("F17").Value = 0.01310043668122271
$ws.Range("H17").Value = 0.1768558951965065
$ws.Range("I17").Value = 0.08951965065502183
$ws.Range("J17").Value = 0.4519650655021834
$ws.Range("K17").Value = 0.08951965065502183
$ws.Range("M17").Value = 0.01091703056768559
$ws.Range("N17").Value = 0.004366812227074236
$ws.Range("O17").Value = 0.04366812227074236
$ws.Range("S17").Value = 0.1200873362445415
$ws.Range("F18").Value = 0.02868852459016394
$ws.Range("H18").Value = 0.2336065573770492
$ws.Range("I18").Value = 0.04918032786885246
$ws.Range("J18").Value = 0.4262295081967213
$ws.Range("K18").Value = 0.0860655737704918
$ws.Range("M18").Value = 0.01639344262295082
$ws.Range("N18").Value = 0.004098360655737705
$ws.Range("O18").Value = 0.05737704918032787
$ws.Range("S18").Value = 0.09836065573770492
$ws.Range("F19").Value = 0.01606732976281561
$ws.Range("H19").Value = 0.2471308339709258
$ws.Range("I19").Value = 0.08339709257842387
$ws.Range("J19").Value = 0.3420045906656465
$ws.Range("K19").Value = 0.1055853098699311
$ws.Range("M19").Value = 0.02065799540933435
$ws.Range("N19").Value = 0.0007651109410864575
$ws.Range("O19").Value = 0.06579954093343535
$ws.Range("S19").Value = 0.1185921958684009

Write-Output "Applied 111 cell updates"
